$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$textStyleRef = $ws.Range("B4")  # plain, unstyled cell used as a style donor to keep General-formatted text cells free of an explicit style index

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.700.39'
$ws.Range("D2").Style = $textStyleRef.Style
$ws.Range("E2").Value = '  -3.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.98'
$ws.Range("D3").Style = $textStyleRef.Style
$ws.Range("E3").Value = '  -3.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = $textStyleRef.Style
$ws.Range("E4").Value = '  -1.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.91'
$ws.Range("D5").Style = $textStyleRef.Style
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("E6").Value = '  -0.99%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4655'
$ws.Range("D7").Style = $textStyleRef.Style
$ws.Range("E7").Value = '  -3.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3913'
$ws.Range("D8").Style = $textStyleRef.Style
$ws.Range("E8").Value = '  -3.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.49'
$ws.Range("D9").Style = $textStyleRef.Style
$ws.Range("E9").Value = '  -2.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07914'
$ws.Range("D10").Style = $textStyleRef.Style
$ws.Range("E10").Value = '  -3.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9856'
$ws.Range("D11").Style = $textStyleRef.Style
$ws.Range("E11").Value = '  -2.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.38'
$ws.Range("D12").Style = $textStyleRef.Style
$ws.Range("E12").Value = '  -5.66%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.852'
$ws.Range("D13").Style = $textStyleRef.Style
$ws.Range("E13").Value = '  -3.70%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.831.19'
$ws.Range("D14").Style = $textStyleRef.Style
$ws.Range("E14").Value = '  -5.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.002'
$ws.Range("D15").Style = $textStyleRef.Style
$ws.Range("E15").Value = '  -3.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06862'
$ws.Range("D16").Style = $textStyleRef.Style
$ws.Range("E16").Value = '  -0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.70'
$ws.Range("D17").Style = $textStyleRef.Style
$ws.Range("E17").Value = '  -4.28%  '

$ws.Range("E19").Value = '  -3.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.13'
$ws.Range("D20").Style = $textStyleRef.Style

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = $textStyleRef.Style
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.699.35'
$ws.Range("D22").Style = $textStyleRef.Style
$ws.Range("E22").Value = '  -3.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.392'
$ws.Range("D23").Style = $textStyleRef.Style
$ws.Range("E23").Value = '  -5.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("D24").Style = $textStyleRef.Style
$ws.Range("E24").Value = '  -5.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.134'
$ws.Range("D25").Style = $textStyleRef.Style
$ws.Range("E25").Value = '  -2.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.042.39'
$ws.Range("D26").Style = $textStyleRef.Style
$ws.Range("E26").Value = '  -5.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.48'
$ws.Range("D27").Style = $textStyleRef.Style
$ws.Range("E27").Value = '  -1.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.50'
$ws.Range("D28").Style = $textStyleRef.Style
$ws.Range("E28").Value = '  -2.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.114'
$ws.Range("D29").Style = $textStyleRef.Style
$ws.Range("E29").Value = '  -5.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.025'
$ws.Range("D30").Style = $textStyleRef.Style
$ws.Range("E30").Value = '  -3.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.64'
$ws.Range("D31").Style = $textStyleRef.Style
$ws.Range("E31").Value = '  -2.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9791'
$ws.Range("D32").Style = $textStyleRef.Style
$ws.Range("E32").Value = '  -3.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09411'
$ws.Range("D33").Style = $textStyleRef.Style
$ws.Range("E33").Value = '  -2.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.374'
$ws.Range("D34").Style = $textStyleRef.Style

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.479'
$ws.Range("D35").Style = $textStyleRef.Style
$ws.Range("E35").Value = '  -2.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.355'
$ws.Range("D36").Style = $textStyleRef.Style
$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("E37").Value = '  -2.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02203'
$ws.Range("D38").Style = $textStyleRef.Style
$ws.Range("E38").Value = '  -3.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.166'
$ws.Range("D39").Style = $textStyleRef.Style
$ws.Range("E39").Value = '  -1.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5730'
$ws.Range("D40").Style = $textStyleRef.Style
$ws.Range("E40").Value = '  -3.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.596'
$ws.Range("D41").Style = $textStyleRef.Style
$ws.Range("E41").Value = '  -3.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.22'
$ws.Range("D42").Style = $textStyleRef.Style
$ws.Range("E42").Value = '  -4.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1803'
$ws.Range("D43").Style = $textStyleRef.Style
$ws.Range("E43").Value = '  -2.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.385'
$ws.Range("D44").Style = $textStyleRef.Style
$ws.Range("E44").Value = '  -2.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.248'
$ws.Range("D45").Style = $textStyleRef.Style
$ws.Range("E45").Value = '  -3.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5400'
$ws.Range("D46").Style = $textStyleRef.Style
$ws.Range("E46").Value = '  -2.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.75'
$ws.Range("D47").Style = $textStyleRef.Style
$ws.Range("E47").Value = '  -5.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07161'
$ws.Range("D48").Style = $textStyleRef.Style
$ws.Range("E48").Value = '  -4.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.915'
$ws.Range("D49").Style = $textStyleRef.Style
$ws.Range("E49").Value = '  -1.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '114.32'
$ws.Range("D50").Style = $textStyleRef.Style
$ws.Range("E50").Value = '  -4.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.28'
$ws.Range("D51").Style = $textStyleRef.Style
$ws.Range("E51").Value = '  +2.95%  '
